# Zeiterfassung_Kasper_Christian.xlsx
# Add two new time-tracking entries (rows 23 and 24) to the "Zeiterfassung" sheet:
#   12.12.2024  5h  Offline-Treffen   Importer, Converter
#   08.01.2025  8h  Fixxes            .obj einlesbar endlich!
# and move the active cell selection to D25 (just below the new data),
# matching the author's commit: "Zeiten aufgeschrieben heute und am 12.12.24"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date number format, alignment, etc.) of the last
# existing data row (22) down into the two new rows so the new cells get
# the same styles (s="4" for the date column, s="2" for the duration column)
# as the rest of the table.
$ws.Range("A22:D22").Copy()
$ws.Range("A23:D24").PasteSpecial(-4122)

# Row 23: 12.12.2024 (serial 45638)
$ws.Cells.Item(23, 1).Value = 45638
$ws.Cells.Item(23, 2).Value = 5
$ws.Cells.Item(23, 3).Value = "Offline-Treffen"
$ws.Cells.Item(23, 4).Value = "Importer, Converter"

# Row 24: 08.01.2025 (serial 45665)
$ws.Cells.Item(24, 1).Value = 45665
$ws.Cells.Item(24, 2).Value = 8
$ws.Cells.Item(24, 3).Value = "Fixxes"
$ws.Cells.Item(24, 4).Value = ".obj einlesbar endlich!"

# Update the selection to reflect where the user ended up after entering
# the new rows.
$null = $ws.Range("D25").Select()
